# Fix property_category label on the "建物" (building) sheet:
# rows 2-7 in column I were mislabeled "land" and should read "building".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "building"
}
